$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = "Death"
$ws.Range("A24").Value = "Jump"
$ws.Range("A26").Value = "Kick"
$ws.Range("A27").Value = "Punch"
$ws.Range("A29").Value = "Jump Attack"

$ws.Range("A27").Select()
